$wb = $excel.ActiveWorkbook

$wsParametry = $wb.Worksheets.Item("Parametry")
$wsSzoki = $wb.Worksheets.Item("Szoki_Rho")

# ---------------------------------------------------------------
# Szoki_Rho (Sheet2): add two new shock-persistence rows (13, 14)
# Formatting copied from the last existing data row (12) so the
# new rows inherit the same styles (A: s7, B: s2, C: s2, D: s5).
# ---------------------------------------------------------------
$wsSzoki.Range("A12:D12").Copy($wsSzoki.Range("A13:D13"))
$wsSzoki.Range("A12:D12").Copy($wsSzoki.Range("A14:D14"))

$wsSzoki.Range("A13").Value = "c_rho_p_energy"
$wsSzoki.Range("A14").Value = "c_rho_p_food"
$wsSzoki.Range("B13").Value = ' ${\rho_p^{energy}}$ '
$wsSzoki.Range("C13").Value = "Persistence of price energy shock"
$wsSzoki.Range("C14").Value = "Persistence of price food shock"
$wsSzoki.Range("B14").Value = '${\rho_p^{food}}$'
$wsSzoki.Range("D13").Value = 0.62406399999999995
$wsSzoki.Range("D14").Value = 0.62406399999999995

# Widen column B now that it holds the longer symbol strings.
$wsSzoki.Columns.Item(2).ColumnWidth = 35.7109375

# ---------------------------------------------------------------
# Parametry (Sheet1): add two new CPI-basket weight rows (26, 27)
# Formatting copied from the last existing data row (25) so the
# new rows inherit the same styles (A: s3, B: s2, C: s2, D: s6).
# ---------------------------------------------------------------
$wsParametry.Range("A25:D25").Copy($wsParametry.Range("A26:D26"))
$wsParametry.Range("A25:D25").Copy($wsParametry.Range("A27:D27"))

$wsParametry.Range("A26").Value = "c_w_energy"
$wsParametry.Range("A27").Value = "c_w_food"
$wsParametry.Range("B27").Value = '${w_{food}}$'
$wsParametry.Range("C27").Value = "Steady state weight - food in CPI basket"
$wsParametry.Range("C26").Value = "Steady state weight - energy in CPI basket"
$wsParametry.Range("B26").Value = '${w_{energy}}$'
$wsParametry.Range("D26").Value = 0.11
$wsParametry.Range("D27").Value = 0.21

# ---------------------------------------------------------------
# View / selection state: Szoki_Rho becomes the active tab, with
# the Parametry sheet scrolled down and showing the new rows
# selected; Szoki_Rho ends up with D14 selected.
# ---------------------------------------------------------------
$wsParametry.Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsParametry.Range("A23:D27").Select()

$wsSzoki.Select()
$wsSzoki.Range("D14").Select()
